$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.493.36'
$ws.Range('E2').Value = '  -1.43%  '

# Row 3
$ws.Range('D3').Value = '1.910.87'
$ws.Range('E3').Value = '  -2.13%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '239.50'
$ws.Range('E5').Value = '  -1.30%  '

# Row 6
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4770'
$ws.Range('E7').Value = '  -2.16%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2844'
$ws.Range('E8').Value = '  -3.31%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06679'
$ws.Range('E9').Value = '  -3.96%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '18.67'
$ws.Range('E10').Value = '  -4.34%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '101.15'
$ws.Range('E11').Value = '  -5.56%  '

# Row 12
$ws.Range('D12').Value = '1.913.87'
$ws.Range('E12').Value = '  -2.09%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07682'
$ws.Range('E13').Value = '  -0.95%  '

# Row 14
$ws.Range('E14').Value = '  -2.06%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6689'
$ws.Range('E15').Value = '  -3.89%  '

# Row 16
$ws.Range('D16').Value = '30.513.43'
$ws.Range('E16').Value = '  -1.43%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '254.11'
$ws.Range('E17').Value = '  -9.09%  '

# Row 18
$ws.Range('E18').Value = '  +0.03%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007469'
$ws.Range('E19').Value = '  -3.87%  '

# Row 20
$ws.Range('E20').Value = '  -4.32%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.374'
$ws.Range('E21').Value = '  -1.81%  '

# Row 22
$ws.Range('E22').Value = '  +0.10%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.4494'
$ws.Range('E23').Value = '  -11.54%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.288'
$ws.Range('E24').Value = '  -3.14%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '168.02'
$ws.Range('E25').Value = '  -0.04%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.329'
$ws.Range('E26').Value = '  -4.15%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.00'
$ws.Range('E27').Value = '  -3.30%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.062'
$ws.Range('E28').Value = '  -5.06%  '

# Row 29
$ws.Range('E29').Value = '  +3.10%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.1007'
$ws.Range('E30').Value = '  -3.28%  '

# Row 31
$ws.Range('E31').Value = '  -0.97%  '

# Row 32
$ws.Range('E32').Value = '  -2.65%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.245'
$ws.Range('E33').Value = '  -3.33%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04711'
$ws.Range('E34').Value = '  -3.29%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7270'
$ws.Range('E35').Value = '  -3.35%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.106'
$ws.Range('E36').Value = '  -4.85%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.000'
$ws.Range('E37').Value = '  +0.12%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.703'
$ws.Range('E38').Value = '  -0.71%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01906'
$ws.Range('E39').Value = '  -4.51%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.613'
$ws.Range('E40').Value = '  -2.41%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '74.85'
$ws.Range('E41').Value = '  -3.26%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.215'
$ws.Range('E42').Value = '  -4.68%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.965'
$ws.Range('E43').Value = '  -6.63%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8609'
$ws.Range('E44').Value = '  -3.99%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '105.04'
$ws.Range('E45').Value = '  -3.39%  '

# Row 46
$ws.Range('B46').Value = 'PaxDollar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.08%  '

# Row 47
$ws.Range('B47').Value = 'TheSandbox'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4238'
$ws.Range('E47').Value = '  -4.26%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.401'
$ws.Range('E48').Value = '  -4.42%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '952.05'
$ws.Range('E49').Value = '  -4.02%  '

# Row 50
$ws.Range('E50').Value = '  -3.88%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '34.75'
$ws.Range('E51').Value = '  -2.90%  '

